$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 5473
$ws1.Range("F4").Value = 11978
$ws1.Range("G4").Value = 58
$ws1.Range("F5").Value = 293
$ws1.Range("F6").Value = 605
$ws1.Range("F7").Value = 177
$ws1.Range("F8").Value = 311
$ws1.Range("F9").Value = 1086
$ws1.Range("F10").Value = 102

# Sheet "全部类型" (All types) updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 5473
$ws4.Range("F7").Value = 11978
$ws4.Range("G7").Value = 58
$ws4.Range("F8").Value = 293
$ws4.Range("F9").Value = 605
$ws4.Range("F10").Value = 177
$ws4.Range("F13").Value = 311
$ws4.Range("F14").Value = 1086
$ws4.Range("F16").Value = 102
